{"js": "// Append two new paragraphs at the very end of the document body:\n//   1) a date line: \"2022\u5e746\u670810\u65e5\u661f\u671f\u4e94\"\n//   2) a journal entry: \"\u591a\u4e91\u8f6c\u5c0f\u96e8\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\u3002\"\n//\n// We build the new content as an OOXML fragment (matching the run/paragraph\n// shape produced by Word itself \u2014 east-Asia font hints on the CJK runs, the\n// digit groups split into their own runs, and a paragraph-mark rPr on the\n// second paragraph) and insert it immediately after the last paragraph in\n// the body, right before the section break.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst insertionPoint = lastParagraph.getRange(\"End\");\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n              <w:t>2</w:t>\n            </w:r>\n            <w:r>\n              <w:t>022</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n              <w:t>\u5e746\u67081</w:t>\n            </w:r>\n            <w:r>\n              <w:t>0</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n              <w:t>\u65e5\u661f\u671f\u4e94</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n              <w:t>\u591a\u4e91\u8f6c\u5c0f\u96e8\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\u3002</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertionPoint.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Append two new paragraphs at the very end of the document body:\n#   1) a date line: \"2022\u5e746\u670810\u65e5\u661f\u671f\u4e94\"\n#   2) a journal entry: \"\u591a\u4e91\u8f6c\u5c0f\u96e8\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\u3002\"\n#\n# Built as an OOXML fragment (matching the run/paragraph shape Word itself\n# produces \u2014 east-Asia font hints on the CJK runs, the digit groups split\n# into their own runs, and a paragraph-mark rPr on the second paragraph)\n# and inserted via Range.InsertXML at the very end of the document, right\n# before the section break.\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Collapse(0)  # wdCollapseEnd\n\n$xml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n              <w:t>2</w:t>\n            </w:r>\n            <w:r>\n              <w:t>022</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n              <w:t>\u5e746\u67081</w:t>\n            </w:r>\n            <w:r>\n              <w:t>0</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n              <w:t>\u65e5\u661f\u671f\u4e94</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>\n              <w:t>\u591a\u4e91\u8f6c\u5c0f\u96e8\uff0c\u4eca\u5929\u5b66\u4e60\u4e86\u5206\u652f\u7ba1\u7406\uff0c\u521b\u5efa\u4e86\u4e00\u4e2adev\u5206\u652f\u3002\u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u53c8\u5feb\u901f\u3002</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$rng.InsertXML($xml)\n"}
